$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D look numeric (e.g. "580.46") but must be stored
# as text, matching the source data (t="inlineStr"/shared-string cells with
# no numeric conversion). Forcing NumberFormat to text before assignment, then
# clearing the format again afterwards, keeps the value as text without leaving
# a stray number-format style on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.996.69'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.381.59'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.85%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.46'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.76'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.619'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +4.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.380.86'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.77%  '

$ws.Range("E10").Value = '  -1.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.91'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.409'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.982.52'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.34%  '

$ws.Range("E14").Value = '  +0.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.91'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.140.42'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("E17").Value = '  +0.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.403.81'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.85'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.71'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.25'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.53'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.75'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.528'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000123'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.75'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.178'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.66%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.98'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.72'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.09'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.98'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.67%  '

$ws.Range("E35").Value = '  -3.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.53'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.47'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.857'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.19'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.39%  '

$ws.Range("E40").Value = '  -0.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.59'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.679.30'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.34'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.23'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0681'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.56'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.63%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.67'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '331.00'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +8.97%  '

$ws.Range("E49").Value = '  -1.41%  '

$ws.Range("E50").Value = '  +2.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.38'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.20%  '
